# Laborator 09.04.2024 - am inceput realizarea unui joc Tower Defence
#
# Mark attendance (TRUE) for "saptamana 7" (column I) for the students
# that were present, on the "Prezente" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$weekPresentRows = 3,5,6,11,18,19,24,29,34,35,36,38,40,46

foreach ($row in $weekPresentRows) {
    $ws.Cells.Item($row, 9).Value = $true
}

# Restore the view/selection state recorded for the sheet: scrolled so
# row 4 is at the top, with I19 as the active/selected cell.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I19").Select()
